$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G (K = strikeouts) with recomputed values
$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 5
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 5
$ws.Range("G7").Value = 3
$ws.Range("G8").Value = 1
